$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "89-82=7"
$t.Cell(1, 2).Range.Text = "58+16=74"
$t.Cell(1, 3).Range.Text = "84-40=44"
$t.Cell(1, 4).Range.Text = "57-46=11"
$t.Cell(1, 5).Range.Text = "78+15=93"
$t.Cell(2, 1).Range.Text = "52-8=44"
$t.Cell(2, 2).Range.Text = "66+6=72"
$t.Cell(2, 3).Range.Text = "16+45=61"
$t.Cell(2, 4).Range.Text = "15+40=55"
$t.Cell(2, 5).Range.Text = "20+60=80"
$t.Cell(3, 1).Range.Text = "87-31=56"
$t.Cell(3, 2).Range.Text = "32-1=31"
$t.Cell(3, 3).Range.Text = "8+89=97"
$t.Cell(3, 4).Range.Text = "72-0=72"
$t.Cell(3, 5).Range.Text = "74+23=97"
$t.Cell(4, 1).Range.Text = "9+14=23"
$t.Cell(4, 2).Range.Text = "46+41=87"
$t.Cell(4, 3).Range.Text = "29+35=64"
$t.Cell(4, 4).Range.Text = "27-12=15"
$t.Cell(4, 5).Range.Text = "33+39=72"
$t.Cell(5, 1).Range.Text = "13+67=80"
$t.Cell(5, 2).Range.Text = "51-11=40"
$t.Cell(5, 3).Range.Text = "80+4=84"
$t.Cell(5, 4).Range.Text = "74+14=88"
$t.Cell(5, 5).Range.Text = "12+16=28"
$t.Cell(6, 1).Range.Text = "41-30=11"
$t.Cell(6, 2).Range.Text = "42-25=17"
$t.Cell(6, 3).Range.Text = "74-48=26"
$t.Cell(6, 4).Range.Text = "84-19=65"
$t.Cell(6, 5).Range.Text = "99-1=98"
$t.Cell(7, 1).Range.Text = "0+33=33"
$t.Cell(7, 2).Range.Text = "82-6=76"
$t.Cell(7, 3).Range.Text = "50-33=17"
$t.Cell(7, 4).Range.Text = "82-28=54"
$t.Cell(7, 5).Range.Text = "37-28=9"
$t.Cell(8, 1).Range.Text = "39-37=2"
$t.Cell(8, 2).Range.Text = "17+67=84"
$t.Cell(8, 3).Range.Text = "47-3=44"
$t.Cell(8, 4).Range.Text = "89-25=64"
$t.Cell(8, 5).Range.Text = "13+76=89"
$t.Cell(9, 1).Range.Text = "86-42=44"
$t.Cell(9, 2).Range.Text = "93-31=62"
$t.Cell(9, 3).Range.Text = "85-78=7"
$t.Cell(9, 4).Range.Text = "0+21=21"
$t.Cell(9, 5).Range.Text = "93-41=52"
$t.Cell(10, 1).Range.Text = "25+69=94"
$t.Cell(10, 2).Range.Text = "38+37=75"
$t.Cell(10, 3).Range.Text = "1-0=1"
$t.Cell(10, 4).Range.Text = "17+62=79"
$t.Cell(10, 5).Range.Text = "88-62=26"
$t.Cell(11, 1).Range.Text = "67-24=43"
$t.Cell(11, 2).Range.Text = "21+17=38"
$t.Cell(11, 3).Range.Text = "23+54=77"
$t.Cell(11, 4).Range.Text = "99-77=22"
$t.Cell(11, 5).Range.Text = "88-8=80"
$t.Cell(12, 1).Range.Text = "31-7=24"
$t.Cell(12, 2).Range.Text = "38-28=10"
$t.Cell(12, 3).Range.Text = "31-22=9"
$t.Cell(12, 4).Range.Text = "93-40=53"
$t.Cell(12, 5).Range.Text = "95-22=73"
$t.Cell(13, 1).Range.Text = "90-1=89"
$t.Cell(13, 2).Range.Text = "17+29=46"
$t.Cell(13, 3).Range.Text = "47-30=17"
$t.Cell(13, 4).Range.Text = "22+18=40"
$t.Cell(13, 5).Range.Text = "33+63=96"
$t.Cell(14, 1).Range.Text = "85-15=70"
$t.Cell(14, 2).Range.Text = "32-0=32"
$t.Cell(14, 3).Range.Text = "38+20=58"
$t.Cell(14, 4).Range.Text = "15+10=25"
$t.Cell(14, 5).Range.Text = "2+11=13"
$t.Cell(15, 1).Range.Text = "75-64=11"
$t.Cell(15, 2).Range.Text = "62-37=25"
$t.Cell(15, 3).Range.Text = "32-22=10"
$t.Cell(15, 4).Range.Text = "39+0=39"
$t.Cell(15, 5).Range.Text = "74-71=3"
$t.Cell(16, 1).Range.Text = "49+29=78"
$t.Cell(16, 2).Range.Text = "87-47=40"
$t.Cell(16, 3).Range.Text = "23-1=22"
$t.Cell(16, 4).Range.Text = "80-7=73"
$t.Cell(16, 5).Range.Text = "82-49=33"
$t.Cell(17, 1).Range.Text = "43-30=13"
$t.Cell(17, 2).Range.Text = "80-19=61"
$t.Cell(17, 3).Range.Text = "34+28=62"
$t.Cell(17, 4).Range.Text = "75-24=51"
$t.Cell(17, 5).Range.Text = "8+82=90"
$t.Cell(18, 1).Range.Text = "33+42=75"
$t.Cell(18, 2).Range.Text = "89-42=47"
$t.Cell(18, 3).Range.Text = "95-56=39"
$t.Cell(18, 4).Range.Text = "84-8=76"
$t.Cell(18, 5).Range.Text = "74-38=36"
$t.Cell(19, 1).Range.Text = "73-66=7"
$t.Cell(19, 2).Range.Text = "69+4=73"
$t.Cell(19, 3).Range.Text = "8+35=43"
$t.Cell(19, 4).Range.Text = "71+0=71"
$t.Cell(19, 5).Range.Text = "9+20=29"
$t.Cell(20, 1).Range.Text = "38-34=4"
$t.Cell(20, 2).Range.Text = "48-18=30"
$t.Cell(20, 3).Range.Text = "73-7=66"
$t.Cell(20, 4).Range.Text = "58+23=81"
$t.Cell(20, 5).Range.Text = "34+23=57"
